{"js": "// Ordered list of (old text -> new text) replacements. Each \"old\" string is\n// unique in the document, so a body-wide literal search for each pair\n// deterministically targets exactly one run: the date heading paragraph\n// plus the 25 populated \"###\u00f7#=##, #\" table cells.\nconst pairs = [\n  [\"2025-09-25 Thursday\", \"2025-09-26 Friday\"],\n  [\"691\u00f75=138, 1\", \"741\u00f73=247, 0\"],\n  [\"463\u00f73=154, 1\", \"245\u00f74=61, 1\"],\n  [\"226\u00f77=32, 2\", \"843\u00f73=281, 0\"],\n  [\"814\u00f76=135, 4\", \"723\u00f78=90, 3\"],\n  [\"769\u00f73=256, 1\", \"146\u00f73=48, 2\"],\n  [\"760\u00f72=380, 0\", \"566\u00f79=62, 8\"],\n  [\"652\u00f72=326, 0\", \"878\u00f73=292, 2\"],\n  [\"778\u00f76=129, 4\", \"998\u00f78=124, 6\"],\n  [\"658\u00f76=109, 4\", \"741\u00f75=148, 1\"],\n  [\"951\u00f73=317, 0\", \"570\u00f73=190, 0\"],\n  [\"139\u00f74=34, 3\", \"455\u00f73=151, 2\"],\n  [\"256\u00f75=51, 1\", \"679\u00f76=113, 1\"],\n  [\"665\u00f76=110, 5\", \"242\u00f73=80, 2\"],\n  [\"475\u00f74=118, 3\", \"377\u00f77=53, 6\"],\n  [\"189\u00f76=31, 3\", \"882\u00f72=441, 0\"],\n  [\"770\u00f76=128, 2\", \"110\u00f79=12, 2\"],\n  [\"780\u00f76=130, 0\", \"481\u00f73=160, 1\"],\n  [\"735\u00f76=122, 3\", \"203\u00f75=40, 3\"],\n  [\"342\u00f77=48, 6\", \"448\u00f74=112, 0\"],\n  [\"226\u00f76=37, 4\", \"960\u00f75=192, 0\"],\n  [\"477\u00f72=238, 1\", \"771\u00f77=110, 1\"],\n  [\"314\u00f77=44, 6\", \"767\u00f72=383, 1\"],\n  [\"398\u00f77=56, 6\", \"984\u00f74=246, 0\"],\n  [\"476\u00f78=59, 4\", \"474\u00f74=118, 2\"],\n  [\"291\u00f76=48, 3\", \"342\u00f78=42, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (old text -> new text) replacements. Each \"old\" string is\n# unique in the document, so a plain Find/Replace (MatchCase, whole literal,\n# no wildcards) for each pair deterministically targets exactly one run:\n# the date heading paragraph plus the 25 populated \"###\u00f7#=##, #\" table cells.\n$pairs = @(\n    @{ Old = \"2025-09-25 Thursday\"; New = \"2025-09-26 Friday\" },\n    @{ Old = \"691\u00f75=138, 1\"; New = \"741\u00f73=247, 0\" },\n    @{ Old = \"463\u00f73=154, 1\"; New = \"245\u00f74=61, 1\" },\n    @{ Old = \"226\u00f77=32, 2\"; New = \"843\u00f73=281, 0\" },\n    @{ Old = \"814\u00f76=135, 4\"; New = \"723\u00f78=90, 3\" },\n    @{ Old = \"769\u00f73=256, 1\"; New = \"146\u00f73=48, 2\" },\n    @{ Old = \"760\u00f72=380, 0\"; New = \"566\u00f79=62, 8\" },\n    @{ Old = \"652\u00f72=326, 0\"; New = \"878\u00f73=292, 2\" },\n    @{ Old = \"778\u00f76=129, 4\"; New = \"998\u00f78=124, 6\" },\n    @{ Old = \"658\u00f76=109, 4\"; New = \"741\u00f75=148, 1\" },\n    @{ Old = \"951\u00f73=317, 0\"; New = \"570\u00f73=190, 0\" },\n    @{ Old = \"139\u00f74=34, 3\"; New = \"455\u00f73=151, 2\" },\n    @{ Old = \"256\u00f75=51, 1\"; New = \"679\u00f76=113, 1\" },\n    @{ Old = \"665\u00f76=110, 5\"; New = \"242\u00f73=80, 2\" },\n    @{ Old = \"475\u00f74=118, 3\"; New = \"377\u00f77=53, 6\" },\n    @{ Old = \"189\u00f76=31, 3\"; New = \"882\u00f72=441, 0\" },\n    @{ Old = \"770\u00f76=128, 2\"; New = \"110\u00f79=12, 2\" },\n    @{ Old = \"780\u00f76=130, 0\"; New = \"481\u00f73=160, 1\" },\n    @{ Old = \"735\u00f76=122, 3\"; New = \"203\u00f75=40, 3\" },\n    @{ Old = \"342\u00f77=48, 6\"; New = \"448\u00f74=112, 0\" },\n    @{ Old = \"226\u00f76=37, 4\"; New = \"960\u00f75=192, 0\" },\n    @{ Old = \"477\u00f72=238, 1\"; New = \"771\u00f77=110, 1\" },\n    @{ Old = \"314\u00f77=44, 6\"; New = \"767\u00f72=383, 1\" },\n    @{ Old = \"398\u00f77=56, 6\"; New = \"984\u00f74=246, 0\" },\n    @{ Old = \"476\u00f78=59, 4\"; New = \"474\u00f74=118, 2\" },\n    @{ Old = \"291\u00f76=48, 3\"; New = \"342\u00f78=42, 6\" }\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($p.Old, $true, $true, $false, $false, $false, $true, 1, $false, $p.New, 2)\n}\n"}
